$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "36.309.93"
$ws.Range("E2").Value = "  +1.82%  "
Set-TextValue $ws.Range("D3") "2.012.01"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "244.54"
$ws.Range("E5").Value = "  -1.03%  "
Set-TextValue $ws.Range("D6") "0.661"
$ws.Range("E6").Value = "  -4.67%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.10%  "
Set-TextValue $ws.Range("D8") "44.42"
$ws.Range("E8").Value = "  +2.97%  "
Set-TextValue $ws.Range("D9") "61.14"
$ws.Range("E9").Value = "  +6.80%  "
$ws.Range("E10").Value = "  +1.18%  "
Set-TextValue $ws.Range("D11") "0.0714"
$ws.Range("E11").Value = "  -5.62%  "
Set-TextValue $ws.Range("D12") "0.0979"
$ws.Range("E12").Value = "  -0.64%  "
Set-TextValue $ws.Range("D13") "14.32"
$ws.Range("E13").Value = "  -1.24%  "
Set-TextValue $ws.Range("D14") "2.299.04"
$ws.Range("E14").Value = "  +5.67%  "
Set-TextValue $ws.Range("D15") "0.802"
$ws.Range("E15").Value = "  -0.28%  "
Set-TextValue $ws.Range("D16") "2.007.34"
$ws.Range("E16").Value = "  +5.66%  "
Set-TextValue $ws.Range("D17") "4.87"
$ws.Range("E17").Value = "  -3.36%  "
Set-TextValue $ws.Range("D18") "36.360.97"
$ws.Range("E18").Value = "  +2.01%  "
Set-TextValue $ws.Range("D19") "71.12"
$ws.Range("E19").Value = "  -3.75%  "
Set-TextValue $ws.Range("D20") "0.0₃0810"
$ws.Range("E20").Value = "  -2.88%  "
Set-TextValue $ws.Range("D21") "236.62"
$ws.Range("E21").Value = "  -4.16%  "
Set-TextValue $ws.Range("D22") "12.70"
$ws.Range("E22").Value = "  -2.43%  "
Set-TextValue $ws.Range("D23") "4.88"
$ws.Range("E23").Value = "  -6.26%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -9.41%  "
Set-TextValue $ws.Range("D26") "165.53"
$ws.Range("E26").Value = "  -0.64%  "
Set-TextValue $ws.Range("D27") "8.61"
$ws.Range("E27").Value = "  -0.53%  "
Set-TextValue $ws.Range("D28") "19.51"
$ws.Range("E28").Value = "  +5.98%  "
Set-TextValue $ws.Range("D29") "1.94"
$ws.Range("E29").Value = "  -9.94%  "
$ws.Range("E30").Value = "  -5.66%  "
Set-TextValue $ws.Range("D31") "21.66"
$ws.Range("E31").Value = "  +48.06%  "
Set-TextValue $ws.Range("D32") "4.32"
$ws.Range("E32").Value = "  -1.40%  "
Set-TextValue $ws.Range("D33") "0.0580"
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D35") "0.0869"
$ws.Range("E35").Value = "  +18.30%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D36") "1.86"
$ws.Range("E36").Value = "  +0.54%  "
Set-TextValue $ws.Range("D37") "3.97"
$ws.Range("E37").Value = "  -6.76%  "
Set-TextValue $ws.Range("D38") "2.13"
$ws.Range("E38").Value = "  +7.76%  "
Set-TextValue $ws.Range("D39") "0.848"
$ws.Range("E39").Value = "  -1.06%  "
Set-TextValue $ws.Range("D40") "1.32"
$ws.Range("E40").Value = "  -11.01%  "
Set-TextValue $ws.Range("D41") "0.0214"
$ws.Range("E41").Value = "  -5.49%  "
Set-TextValue $ws.Range("D42") "95.32"
$ws.Range("E42").Value = "  -4.09%  "
Set-TextValue $ws.Range("D43") "1.10"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("E44").Value = "  +14.87%  "
Set-TextValue $ws.Range("D45") "15.80"
$ws.Range("E45").Value = "  -7.76%  "
Set-TextValue $ws.Range("D46") "1.303.72"
$ws.Range("E46").Value = "  -1.48%  "
Set-TextValue $ws.Range("D47") "0.0813"
$ws.Range("E47").Value = "  -0.05%  "
Set-TextValue $ws.Range("D48") "2.77"
$ws.Range("E48").Value = "  +0.86%  "
Set-TextValue $ws.Range("D49") "2.192.08"
$ws.Range("E49").Value = "  +5.48%  "
Set-TextValue $ws.Range("D50") "2.19"
$ws.Range("E50").Value = "  -7.89%  "
Set-TextValue $ws.Range("D51") "3.82"
$ws.Range("E51").Value = "  +15.48%  "
